$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.399.48'
$ws.Range('E2').Value = '  +3.45%  '
$ws.Range('D3').Value = '2.009.75'
$ws.Range('E3').Value = '  +7.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7627'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +61.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '258.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.86%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9990'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3602'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +25.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '29.01'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +33.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07087'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8554'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +19.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08121'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '102.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.21%  '
$ws.Range('D14').Value = '2.011.56'
$ws.Range('E14').Value = '  +7.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.622'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +9.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '273.98'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.36%  '
$ws.Range('D17').Value = '31.399.14'
$ws.Range('E17').Value = '  +3.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +12.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.943'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +12.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000008013'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.94%  '
$ws.Range('D21').Value = '2.272.32'
$ws.Range('E21').Value = '  +7.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9991'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9989'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.331'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +17.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.98'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1461'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +51.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.361'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +25.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.616'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.649'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.354'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.411'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05221'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.238'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +10.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7646'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +11.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.797'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02024'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.946'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.773'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '80.96'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.193'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +14.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4783'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +14.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8634'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.99'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.708'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.01'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4399'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +12.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '951.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.89%  '
